# chore(auto): pipeline update & logic validation
#
# "super" (row 38) drops out of the war-participation list; every
# subsequent participant (RobaFrag..Alvaro, rows 39-48) moves up one row,
# and "super" is re-appended at the bottom of that block (row 48) with
# their original status/scores intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the 11-row block (A38:H48) exactly as it stands before the shift.
$firstRow = 38
$lastRow  = 48
$numRows  = $lastRow - $firstRow + 1

$names    = @()
$statuses = @()
$sources  = @()
$scores   = @()

for ($i = 0; $i -lt $numRows; $i++) {
    $r = $firstRow + $i
    $names    += $ws.Cells.Item($r, 1).Value2
    $statuses += $ws.Cells.Item($r, 2).Value2
    $sources  += $ws.Cells.Item($r, 3).Value2
    $rowScores = @()
    for ($c = 4; $c -le 8; $c++) {
        $rowScores += $ws.Cells.Item($r, $c).Value2
    }
    $scores += ,$rowScores
}

# Row 38 ("super") is removed from the top; rows 39..48 shift up into
# 38..47, and "super" (originally at index 0) is appended at row 48.
for ($i = 0; $i -lt ($numRows - 1); $i++) {
    $destRow = $firstRow + $i
    $srcIdx  = $i + 1

    $ws.Cells.Item($destRow, 1).Value = $names[$srcIdx]
    $ws.Cells.Item($destRow, 2).Value = $statuses[$srcIdx]
    $ws.Cells.Item($destRow, 3).Value = $sources[$srcIdx]
    for ($c = 4; $c -le 8; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $scores[$srcIdx][$c - 4]
    }
}

# Re-append "super" (original index 0) at the bottom of the block, row 48.
$ws.Cells.Item($lastRow, 1).Value = $names[0]
$ws.Cells.Item($lastRow, 2).Value = $statuses[0]
$ws.Cells.Item($lastRow, 3).Value = $sources[0]
for ($c = 4; $c -le 8; $c++) {
    $ws.Cells.Item($lastRow, $c).Value = $scores[0][$c - 4]
}
